$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10369
$ws1.Range("F3").Value = 429
$ws1.Range("F4").Value = 2530
$ws1.Range("F9").Value = 777
$ws1.Range("F11").Value = 1219
$ws1.Range("F13").Value = 3234
$ws1.Range("F14").Value = 2407
$ws1.Range("F16").Value = 2158
$ws1.Range("F17").Value = 2158
$ws1.Range("F20").Value = 485
$ws1.Range("F23").Value = 67
$ws1.Range("F28").Value = 50
$ws1.Range("F29").Value = 378
$ws1.Range("F32").Value = 393
$ws1.Range("F33").Value = 603
$ws1.Range("F34").Value = 16
$ws1.Range("E35").Value = "2024.04.30 10:00-05.01 16:00"
$ws1.Range("F36").Value = 260
$ws1.Range("F38").Value = 1579
$ws1.Range("F39").Value = 468
$ws1.Range("F40").Value = 456
$ws1.Range("F41").Value = 1716
$ws1.Range("F42").Value = 139
$ws1.Range("F43").Value = 444
$ws1.Range("F44").Value = 54
$ws1.Range("F45").Value = 459
$ws1.Range("F46").Value = 1026
$ws1.Range("F48").Value = 365

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value = 266

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10369
$ws4.Range("F3").Value = 429
$ws4.Range("F4").Value = 2530
$ws4.Range("G7").Value = 266
$ws4.Range("F11").Value = 777
$ws4.Range("F13").Value = 3234
$ws4.Range("F14").Value = 2407
$ws4.Range("F15").Value = 2158
$ws4.Range("F16").Value = 2158
$ws4.Range("F19").Value = 67
$ws4.Range("F24").Value = 50
$ws4.Range("F25").Value = 378
$ws4.Range("F28").Value = 393
$ws4.Range("F29").Value = 603
$ws4.Range("F30").Value = 16
$ws4.Range("E34").Value = "2024.04.30 10:00-05.01 16:00"
$ws4.Range("F35").Value = 260
$ws4.Range("F36").Value = 1579
$ws4.Range("F37").Value = 468
$ws4.Range("F39").Value = 456
$ws4.Range("F40").Value = 1716
$ws4.Range("F41").Value = 139
$ws4.Range("F45").Value = 444
$ws4.Range("F46").Value = 54
$ws4.Range("F47").Value = 459
$ws4.Range("F48").Value = 1026
$ws4.Range("F49").Value = 365
